# Meeting Stuff & RAS Hotfix
#
# Logs actual hours burned on "Day 4" (column F) for three tasks that were
# missed in the original entry, and leaves the selection where the edit
# was made. The downstream "Actual Burndown" cells (G3/H3), which are
# formulas (=$C3-SUM(C7:F38) / =$C3-SUM(C7:G38)), recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Upper Floor Room F - 2 hours on Day 4
$ws.Range("F10").Value = 2

# Upper Floor Room I - 2 hours on Day 4
$ws.Range("F11").Value = 2

# Gamepad Input - 1 hour on Day 4
$ws.Range("F19").Value = 1

# Leave the cursor on the last-touched cell, matching the saved selection.
$ws.Range("G15").Select()
